$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicated header row (row 32). Everything below shifts up by one row.
$ws.Rows("32:32").Delete()

# Update the view to match the post-edit state: scrolled so row 24 is at top,
# with row 32 (now the first data row of the "OpenAI 4o" block) selected.
$ws.Range("A32:XFD32").Select()
$excel.ActiveWindow.ScrollRow = 24
